$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BOM content changes (asm/Pick_Place.xlsx "Pick_Place" sheet) ---

# 1. C6's Comment: 30pF -> 33pF
$rowC6 = $ws.Columns(1).Find("C6").Row
$ws.Cells.Item($rowC6, 2).Value = "33pF"
# Re-apply the original cell formatting (border/font/quote-style) that a plain
# value write resets, by pasting formats from the untouched cell above it.
$ws.Cells.Item($rowC6 - 1, 2).Copy()
$ws.Cells.Item($rowC6, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. U2's Comment: NCP702N -> TPS78233
$rowU2 = $ws.Columns(1).Find("U2").Row
$ws.Cells.Item($rowU2, 2).Value = "TPS78233"
$ws.Cells.Item($rowU2 - 1, 2).Copy()
$ws.Cells.Item($rowU2, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Add a new row for designator R9 (10K / 0603_res), placed right before R10
#    (matches the BOM gap: R8 was directly followed by R10 before this edit)
$rowR10 = $ws.Columns(1).Find("R10").Row
$ws.Rows($rowR10).Insert()
$newRow = $rowR10
$ws.Cells.Item($newRow, 1).Value = "R9"
$ws.Cells.Item($newRow, 2).Value = "10K"
$ws.Cells.Item($newRow, 3).Value = "0603_res"
# Inherit the row formatting (borders/font/style) from the row below it, since a
# freshly-inserted row otherwise comes back with the default (unformatted) style.
$ws.Range("A" + ($newRow + 1) + ":C" + ($newRow + 1)).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. Remove the row for designator R13 (component no longer populated)
$rowR13 = $ws.Columns(1).Find("R13").Row
$ws.Rows($rowR13).Delete()
